$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pool_seq_analyses")

$ws.Range("E5").Value = "Eddie 49330882"
$ws.Range("F5").Value = 45830
$ws.Range("F5").NumberFormat = $ws.Range("D5").NumberFormat
$ws.Range("G5").Value = "NA"

$ws.Range("A6").Value = "Set_9_std"
$ws.Range("B6").Value = "Standard set with burnin analysed with simulations pf poolseq (read_length = 75, coverage = 1000, V_logmean = 0)"
$ws.Range("C6").Value = 100
$ws.Range("D6").Value = 45830
$ws.Range("D6").NumberFormat = $ws.Range("D5").NumberFormat

$ws.Range("F5:G5").Select()
